{"js": "// Insert a new worklog entry paragraph right after the paragraph that\n// ends with \"...get_apod_info function.\" and before the trailing blank\n// paragraphs, matching the author's new 2022-04-26 entry.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text ends with \"function.\" (the existing\n// 2022-04-24 entry) so the new entry is inserted directly after it.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  if (text.trim().endsWith(\"function.\")) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the '...function.' paragraph to insert after.\");\n}\n\n// Collapse to a range immediately after that paragraph, then insert the\n// new paragraph's OOXML there so we can control the run/proofErr\n// boundaries exactly (mirrors how Word itself records spell-check\n// boundaries around the underscore-separated identifiers).\nconst insertionPoint = target.getRange(Word.RangeLocation.after);\nawait context.sync();\n\nconst newParagraphOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n<w:r><w:t xml:space=\"preserve\">2022-04-26 created the </w:t></w:r>\n<w:r><w:t>create _</w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>image_db</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> and </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>add_image_to_db</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> functions. Used the techniques from lab</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> 3 to remind myself how to work with databases.</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:t>This took me about 1 hour.</w:t></w:r>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\n\ninsertionPoint.insertOoxml(newParagraphOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert a new worklog entry paragraph right after the paragraph that\n# ends with \"...get_apod_info function.\" and before the trailing blank\n# paragraphs, matching the author's new 2022-04-26 entry.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text ends with \"function.\" (the existing\n# 2022-04-24 entry) so the new entry is inserted directly after it.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\", \" \")\n    if ($t.EndsWith(\"function.\")) {\n        $target = $p\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the '...function.' paragraph to insert after.\"\n}\n\n# Collapse to a point immediately after that paragraph (its trailing\n# paragraph mark), then insert the new paragraph's OOXML there so we can\n# control the run/proofErr boundaries exactly (mirrors how Word itself\n# records spell-check boundaries around the underscore-separated\n# identifiers).\n$insertionPoint = $target.Range.Duplicate\n$insertionPoint.Collapse(0)  # wdCollapseEnd\n\n$newParagraphOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n<w:r><w:t xml:space=\"preserve\">2022-04-26 created the </w:t></w:r>\n<w:r><w:t>create _</w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>image_db</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> and </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>add_image_to_db</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> functions. Used the techniques from lab</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> 3 to remind myself how to work with databases.</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:t>This took me about 1 hour.</w:t></w:r>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n\"@\n\n$insertionPoint.InsertXML($newParagraphOoxml, \"After\")\n"}
